$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split "?      Evet ( " into "?      Evet (" + " " (two bold runs).
#    We locate the run via Find, then force Word to split off the final
#    trailing space character into its own run by re-asserting identical
#    (already-true) Bold formatting just on that character.
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("?      Evet ( ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find '?      Evet ( '" }
$spaceRange = $d.Range($r.End - 1, $r.End)
$spaceRange.Bold = 0
$spaceRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Change the hidden MERGEFIELD instruction from "Y1" to "y1", then
#    force Word to refresh/merge the displayed field result run(s)
#    ("«" + "y" + "1»" -> single "«y1»" run).
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $fld = $d.Fields.Item($i)
    if ($fld.Code.Text -eq " MERGEFIELD  Y1  \* MERGEFORMAT ") {
        $target = $fld
        break
    }
}
if ($null -eq $target) { throw "Could not find MERGEFIELD Y1" }
$target.Code = " MERGEFIELD  y1  \* MERGEFORMAT "

# Merge the (now stale, still 3-run) field-result display "«" + "y" + "1»"
# into a single run by replacing the result range's text with itself via
# Find/Replace (format-preserving) -- this normalizes the backing runs
# without invoking a document-wide Fields.Update(), which would otherwise
# also disturb unrelated MERGEFIELD results elsewhere in the document.
$resultRange = $target.Result
$resultText = $resultRange.Text
$null = $resultRange.Find.Execute($resultText, $false, $false, $false, $false, $false, $true, 1, $false, $resultText, 2)

# ---------------------------------------------------------------------
# 3) Split " )     Hayır (" into " " + ")     Hayır (" (two bold runs).
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(" )     Hayır (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find ' )     Hayır ('" }
$leadSpaceRange = $d.Range($r2.Start, $r2.Start + 1)
$leadSpaceRange.Bold = 0
$leadSpaceRange.Bold = 1
